$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 34: Rainette crucifère - Jonquière, Saguenay - Martin Bertrand
# ---------------------------------------------------------------------------
$ws.Range("A34").Value = 45776
$ws.Range("B34").Value = "Rainette crucifère"
$ws.Range("C34").Value = "N/A"
$ws.Range("D34").Value = "Jonquière, Saguenay"
$ws.Range("E34").Value = "Saguenay–Lac-Saint-Jean"
$ws.Range("F34").Value = "D"
$ws.Range("G34").Value = "Cote 1"
$ws.Range("H34").Value = "Un seul individu ; il est encore trop tôt pour commencer le premier inventaire."
$ws.Range("I34").Value = "Martin Bertrand"

# ---------------------------------------------------------------------------
# Row 35: Grenouille des bois - Saguenay - Martin Bertrand
# ---------------------------------------------------------------------------
$ws.Range("A35").Value = 45778
$ws.Range("B35").Value = "Grenouille des bois"
$ws.Range("C35").Value = 66
$ws.Range("D35").Value = "Saguenay"
$ws.Range("E35").Value = "Saguenay–Lac-Saint-Jean"
$ws.Range("F35").Value = "C"
$ws.Range("G35").Value = "Cote 3"
$ws.Range("I35").Value = "Martin Bertrand"

# ---------------------------------------------------------------------------
# Row 36: Rainette crucifère - Lac-Ministuk - Yoann Perrot
# ---------------------------------------------------------------------------
$ws.Range("A36").Value = 45778
$ws.Range("B36").Value = "Rainette crucifère"
$ws.Range("C36").Value = "201m"
$ws.Range("D36").Value = "Lac-Ministuk"
$ws.Range("E36").Value = "Saguenay–Lac-Saint-Jean"
$ws.Range("F36").Value = "C"
$ws.Range("G36").Value = "Cote 3"
$ws.Range("I36").Value = "Yoann Perrot"

# ---------------------------------------------------------------------------
# Row 37: Grenouille des bois - Lac-Ministuk - Yoann Perrot
# ---------------------------------------------------------------------------
$ws.Range("A37").Value = 45778
$ws.Range("B37").Value = "Grenouille des bois"
$ws.Range("C37").Value = "201m"
$ws.Range("D37").Value = "Lac-Ministuk"
$ws.Range("E37").Value = "Saguenay–Lac-Saint-Jean"
$ws.Range("F37").Value = "C"
$ws.Range("G37").Value = "Cote 1"
$ws.Range("I37").Value = "Yoann Perrot"

# ---------------------------------------------------------------------------
# Formatting: reuse existing cell styles by copying formats from cells that
# already carry them (date format, species fill, plain/centered cells, etc.)
# ---------------------------------------------------------------------------

# Date cells (same format as the rest of column A)
$ws.Range("A4").Copy()
$ws.Range("A34:A37").PasteSpecial(-4122)

# "Rainette crucifère" rows use the same fill as other male/"s=9" entries
$ws.Range("B6").Copy()
$ws.Range("B34").PasteSpecial(-4122)
$ws.Range("B36").PasteSpecial(-4122)

# "Grenouille des bois" rows use the same fill as other female/"s=10" entries
$ws.Range("B5").Copy()
$ws.Range("B35").PasteSpecial(-4122)
$ws.Range("B37").PasteSpecial(-4122)

# Generic centered cells (C, D, G, I columns)
$ws.Range("J3").Copy()
$ws.Range("C34:D34").PasteSpecial(-4122)
$ws.Range("G34").PasteSpecial(-4122)
$ws.Range("I34").PasteSpecial(-4122)
$ws.Range("C35:D35").PasteSpecial(-4122)
$ws.Range("G35").PasteSpecial(-4122)
$ws.Range("I35").PasteSpecial(-4122)
$ws.Range("C36:D36").PasteSpecial(-4122)
$ws.Range("G36").PasteSpecial(-4122)
$ws.Range("I36").PasteSpecial(-4122)
$ws.Range("C37:D37").PasteSpecial(-4122)
$ws.Range("G37").PasteSpecial(-4122)
$ws.Range("I37").PasteSpecial(-4122)

# Zone climatique "C" cells reuse the existing yellow fill style
$ws.Range("F28").Copy()
$ws.Range("F35").PasteSpecial(-4122)
$ws.Range("F36").PasteSpecial(-4122)
$ws.Range("F37").PasteSpecial(-4122)

# Zone climatique "D" is a brand-new category: green fill, centered text
$ws.Range("F34").Interior.Color = 5296274
$ws.Range("F34").HorizontalAlignment = -4108

$ws.Range("A1").Select()
